$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 151.7260716666667
$ws.Range("H2").Value = 455.178215
$ws.Range("I2").Value = 0.2700739458961593
$ws.Range("J2").Value = 0.2783366498663096
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.743532666666667
$ws.Range("N2").Value = 11.230598
$ws.Range("O2").Value = 0.9802973346235675
$ws.Range("P2").Value = 0.9802973346235675
$ws.Range("Q2").Value = 567.9915056691746
$ws.Range("R2").Value = 5111.923551022571
$ws.Range("S2").Value = 0.2647527693132745
$ws.Range("T2").Value = 0.2728526759919964

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 151.7260716666667
$ws.Range("H3").Value = 455.178215
$ws.Range("I3").Value = 0.2700739458961593
$ws.Range("J3").Value = 0.2783366498663096
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.07524
$ws.Range("N3").Value = 0.22572
$ws.Range("O3").Value = 0.01970266537643246
$ws.Range("P3").Value = 0.01970266537643246
$ws.Range("Q3").Value = 11.4158696322
$ws.Range("R3").Value = 102.7428266898
$ws.Range("S3").Value = 0.005321176582884752
$ws.Range("T3").Value = 0.005483973874313141

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 82.248871
$ws.Range("H4").Value = 246.746613
$ws.Range("I4").Value = 0.146403824289839
$ws.Range("J4").Value = 0.150882936320401
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.743532666666667
$ws.Range("N4").Value = 11.230598
$ws.Range("O4").Value = 0.9802973346235675
$ws.Range("P4").Value = 0.9802973346235675
$ws.Range("Q4").Value = 307.9013353849526
$ws.Range("R4").Value = 2771.112018464574
$ws.Range("S4").Value = 0.1435192787300263
$ws.Range("T4").Value = 0.1479101403150666

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 82.248871
$ws.Range("H5").Value = 246.746613
$ws.Range("I5").Value = 0.146403824289839
$ws.Range("J5").Value = 0.150882936320401
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07524
$ws.Range("N5").Value = 0.22572
$ws.Range("O5").Value = 0.01970266537643246
$ws.Range("P5").Value = 0.01970266537643246
$ws.Range("Q5").Value = 6.18840505404
$ws.Range("R5").Value = 55.69564548636
$ws.Range("S5").Value = 0.002884545559812712
$ws.Range("T5").Value = 0.002972796005334428

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 123.444321
$ws.Range("H6").Value = 370.332963
$ws.Range("I6").Value = 0.2197321429647646
$ws.Range("J6").Value = 0.2264546783208506
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.743532666666667
$ws.Range("N6").Value = 11.230598
$ws.Range("O6").Value = 0.9802973346235675
$ws.Range("P6").Value = 0.9802973346235675
$ws.Range("Q6").Value = 462.117848177986
$ws.Range("R6").Value = 4159.060633601875
$ws.Range("S6").Value = 0.2154028340794834
$ws.Range("T6").Value = 0.2219929175709672

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Agtr2"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 123.444321
$ws.Range("H7").Value = 370.332963
$ws.Range("I7").Value = 0.2197321429647646
$ws.Range("J7").Value = 0.2264546783208506
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.07524
$ws.Range("N7").Value = 0.22572
$ws.Range("O7").Value = 0.01970266537643246
$ws.Range("P7").Value = 0.01970266537643246
$ws.Range("Q7").Value = 9.28795071204
$ws.Range("R7").Value = 83.59155640836
$ws.Range("S7").Value = 0.004329308885281175
$ws.Range("T7").Value = 0.004461760749883373

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Agtr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 154.3429766666667
$ws.Range("H8").Value = 463.02893
$ws.Range("I8").Value = 0.2747320633285943
$ws.Range("J8").Value = 0.2831372788071194
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.743532666666667
$ws.Range("N8").Value = 11.230598
$ws.Range("O8").Value = 0.9802973346235675
$ws.Range("P8").Value = 0.9802973346235675
$ws.Range("Q8").Value = 577.7879750222378
$ws.Range("R8").Value = 5200.091775200141
$ws.Range("S8").Value = 0.2693191094166542
$ws.Range("T8").Value = 0.2775587197471891

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Agtr2"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 154.3429766666667
$ws.Range("H9").Value = 463.02893
$ws.Range("I9").Value = 0.2747320633285943
$ws.Range("J9").Value = 0.2831372788071194
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.07524
$ws.Range("N9").Value = 0.22572
$ws.Range("O9").Value = 0.01970266537643246
$ws.Range("P9").Value = 0.01970266537643246
$ws.Range("Q9").Value = 11.6127655644
$ws.Range("R9").Value = 104.5148900796
$ws.Range("S9").Value = 0.005412953911940146
$ws.Range("T9").Value = 0.005578559059930335

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Agtr2"
$ws.Range("D10").Value = "FAPs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 50.0323125
$ws.Range("H10").Value = 100.064625
$ws.Range("I10").Value = 0.08905802352064279
$ws.Range("J10").Value = 0.06118845668531954
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.743532666666667
$ws.Range("N10").Value = 11.230598
$ws.Range("O10").Value = 0.9802973346235675
$ws.Range("P10").Value = 0.9802973346235675
$ws.Range("Q10").Value = 187.297596232625
$ws.Range("R10").Value = 1123.78557739575
$ws.Range("S10").Value = 0.08730334308412911
$ws.Range("T10").Value = 0.05998288099834835

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Agtr2"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 50.0323125
$ws.Range("H11").Value = 100.064625
$ws.Range("I11").Value = 0.08905802352064279
$ws.Range("J11").Value = 0.06118845668531954
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.07524
$ws.Range("N11").Value = 0.22572
$ws.Range("O11").Value = 0.01970266537643246
$ws.Range("P11").Value = 0.01970266537643246
$ws.Range("Q11").Value = 3.7644311925
$ws.Range("R11").Value = 22.586587155
$ws.Range("S11").Value = 0.001754680436513677
$ws.Range("T11").Value = 0.001205575686971183
